$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT without letting Excel's "smart"
# type detection silently turn numeric-looking strings (e.g. "1",
# "17.0000", "0:4") into real numbers / other types. Temporarily forcing
# the format to Text ("@") for the assignment - then restoring the
# cell's original number format - keeps the cell's existing style
# (these report cells already carry a numeric-looking display format
# even though they hold literal text).
function Set-TextValue($range, [string]$value) {
    $orig = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $orig
}

# The report grew from 2 line items to 4 line items:
#   - a new item ("DEPOVIT ...") is inserted as the new #1 (current row 7)
#   - the previous #1 (PANADOL ...) becomes #2 (current row 8)
#   - the previous #2 (POWER B COMPLEX ...) becomes #3 (new row 9)
#   - a new item ("صوفى طويل جدا جدا") becomes #4 (new row 10)
#   - the totals row moves from row 9 -> row 11 (with an updated sum)
#   - the footer row moves from row 10 -> row 12 (with an updated timestamp)
#
# Capture the two rows that will be pushed down (the old totals row and the
# old footer row) before anything is overwritten, then duplicate the
# formatting of rows 7/8 into the new item rows 9/10.

# 1) Push the current totals row (9) and footer row (10) down to 11 / 12,
#    carrying their formatting (styles, merges) with them.
$ws.Range("A10:Q10").Copy($ws.Range("A12:Q12"))
$ws.Range("A9:Q9").Copy($ws.Range("A11:Q11"))

# 2) Duplicate the formatting of the existing item rows (7, 8) into the new
#    item rows (9, 10) that will hold items #3 and #4.
$ws.Range("A7:Q7").Copy($ws.Range("A9:Q9"))
$ws.Range("A8:Q8").Copy($ws.Range("A10:Q10"))

# 3) Restore the correct row heights (Copy() does not carry row height).
$ws.Rows(9).RowHeight = 25.5
$ws.Rows(10).RowHeight = 24.75
$ws.Rows(11).RowHeight = 25.5
$ws.Rows(12).RowHeight = 16.5

# 4) Move the (old #2) POWER B COMPLEX line into its new spot, row 9, item #3.
$ws.Range("A9").Value = 3
Set-TextValue $ws.Range("C9") "POWER B COMPLEX I.M./I.V. 6 AMP"
Set-TextValue $ws.Range("H9") "0:3"
Set-TextValue $ws.Range("L9") "1"
Set-TextValue $ws.Range("N9") "48.00"
Set-TextValue $ws.Range("P9") "-7.6800"
Set-TextValue $ws.Range("Q9") "0:-1"

# 5) Write the brand-new item #4 into row 10.
$ws.Range("A10").Value = 4
Set-TextValue $ws.Range("C10") "صوفى طويل جدا جدا"
Set-TextValue $ws.Range("H10") "12:0"
Set-TextValue $ws.Range("L10") "0"
Set-TextValue $ws.Range("N10") "55.00"
Set-TextValue $ws.Range("P10") "55.0000"
Set-TextValue $ws.Range("Q10") "1:0"

# 6) Move the (old #1) PANADOL line down into row 8, item #2.
Set-TextValue $ws.Range("C8") "PANADOL ADVANCE 500 MG 48 TABLETS"
Set-TextValue $ws.Range("H8") "1:2"
Set-TextValue $ws.Range("L8") "1"
Set-TextValue $ws.Range("N8") "92.00"
Set-TextValue $ws.Range("P8") "46.0000"
Set-TextValue $ws.Range("Q8") "0:2"
$ws.Range("A8").Value = 2

# 7) Write the brand-new item #1 into row 7.
Set-TextValue $ws.Range("C7") "DEPOVIT B12-1000MCG/ML 5 I.M. AMP"
Set-TextValue $ws.Range("H7") "0:4"
Set-TextValue $ws.Range("L7") "1"
Set-TextValue $ws.Range("N7") "85.00"
Set-TextValue $ws.Range("P7") "17.0000"
Set-TextValue $ws.Range("Q7") "0:1"
$ws.Range("A7").Value = 1

# 8) Update the totals row (now row 11) with the new sum of the sell prices.
$ws.Range("P11").Value = 110.31999999999999

# 9) Update the footer row (now row 12) with the new generation timestamp.
Set-TextValue $ws.Range("A12") "Thursday, 7 August, 2025 9:58 AM"

Write-Output "done"
